$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    # Force the assignment to be treated as plain text so Excel does not
    # auto-convert the dd-mm-yyyy looking string into a date serial number.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    # Restore the cell's original (default) style now that the text value
    # is safely stored, so no stray style index is left on the cell.
    $cell.Style = "Normal"
}
